$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are treated as text (not auto-converted to numbers)
# while keeping the original (default/Normal) cell style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.996.22"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "2.779.47"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "358.13"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").Value = "109.37"
$ws.Range("E6").Value = "  -4.23%  "
$ws.Range("D7").Value = "0.564"
$ws.Range("E7").Value = "  +3.07%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("D10").Value = "40.12"
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "0.132"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "19.44"
$ws.Range("E13").Value = "  -3.44%  "
$ws.Range("D14").Value = "7.61"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "3.215.61"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").Value = "2.791.58"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").Value = "0.923"
$ws.Range("E17").Value = "  +3.59%  "
$ws.Range("D18").Value = "51.870.88"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "7.38"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "13.07"
$ws.Range("E21").Value = "  -5.40%  "
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").Value = "273.93"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").Value = "69.63"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").Value = "26.57"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "10.13"
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").Value = "0.144"
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("D31").Value = "0.0463"
$ws.Range("E31").Value = "  +4.09%  "
$ws.Range("D32").Value = "51.43"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").Value = "34.10"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").Value = "5.72"
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("D35").Value = "5.44"
$ws.Range("E35").Value = "  +10.80%  "
$ws.Range("D36").Value = "0.0837"
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "3.21"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").Value = "18.27"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E40").Value = "  -4.38%  "
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").Value = "124.07"
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("E44").Value = "  -2.22%  "
$ws.Range("D45").Value = "21.81"
$ws.Range("E45").Value = "  -7.27%  "
$ws.Range("D46").Value = "2.066.61"
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("D47").Value = "3.25"
$ws.Range("E47").Value = "  -3.56%  "
$ws.Range("D48").Value = "2.28"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "0.926"
$ws.Range("E50").Value = "  -4.67%  "
$ws.Range("D51").Value = "8.96"
$ws.Range("E51").Value = "  +0.70%  "

# Restore default styling on the Price column so no residual text format remains applied.
$ws.Range("D2:D51").Style = "Normal"
